# Updates the cryptocurrency price/volume table (columns D and E) in place,
# forcing values to remain text (matching the source feed formatting) even
# when they look numeric, and resetting the style afterwards so no visible
# formatting change is introduced.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2" = "41.655.62"
    "E2" = "  +0.14%  "
    "D3" = "2.475.17"
    "E3" = "  +0.68%  "
    "E4" = "  +0.06%  "
    "D5" = "319.37"
    "E5" = "  +1.51%  "
    "D6" = "92.37"
    "E6" = "  +0.22%  "
    "E7" = "  +0.83%  "
    "D9" = "0.511"
    "D10" = "0.0861"
    "E10" = "  +8.09%  "
    "D11" = "33.06"
    "E11" = "  +2.23%  "
    "E12" = "  -0.07%  "
    "D13" = "2.857.35"
    "E13" = "  +0.76%  "
    "D14" = "6.88"
    "E14" = "  +0.64%  "
    "D15" = "15.51"
    "E15" = "  -1.71%  "
    "D16" = "2.482.45"
    "E16" = "  +1.54%  "
    "D17" = "0.793"
    "E17" = "  +2.28%  "
    "D18" = "41.601.23"
    "E18" = "  +0.01%  "
    "E19" = "  -0.29%  "
    "D20" = "0.0₃0941"
    "E20" = "  +0.65%  "
    "D21" = "70.65"
    "E21" = "  -0.04%  "
    "D22" = "11.22"
    "E22" = "  -0.64%  "
    "D23" = "239.73"
    "E23" = "  +0.81%  "
    "E24" = "  +1.65%  "
    "E25" = "  +2.29%  "
    "D27" = "24.95"
    "E27" = "  +2.82%  "
    "D28" = "2.24"
    "E28" = "  -0.51%  "
    "E29" = "  +0.58%  "
    "D30" = "36.65"
    "E30" = "  +4.58%  "
    "D31" = "157.00"
    "E31" = "  +0.73%  "
    "D32" = "5.42"
    "E32" = "  -0.43%  "
    "E33" = "  +0.03%  "
    "D34" = "0.0763"
    "E34" = "  +0.68%  "
    "E35" = "  -0.26%  "
    "D36" = "17.18"
    "E36" = "  -1.64%  "
    "D37" = "0.116"
    "E37" = "  +1.82%  "
    "E38" = "  +2.55%  "
    "D39" = "2.88"
    "E39" = "  +0.10%  "
    "E40" = "  +1.48%  "
    "D41" = "4.00"
    "E41" = "  +1.72%  "
    "D42" = "2.46"
    "E42" = "  -1.33%  "
    "D43" = "2.001.60"
    "E43" = "  +1.59%  "
    "E44" = "  +0.71%  "
    "D45" = "18.67"
    "E45" = "  -1.30%  "
    "D46" = "2.97"
    "E46" = "  +2.25%  "
    "E47" = "  +6.34%  "
    "D48" = "2.718.29"
    "E48" = "  +0.89%  "
    "D49" = "97.58"
    "E49" = "  +0.89%  "
    "D50" = "75.59"
    "E50" = "  +5.42%  "
    "D51" = "67.23"
    "E51" = "  +1.02%  "
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.Style = "Normal"
}
